# Auto-generated: applies the cryptos.xlsx price/volume/coin-identity refresh
# described by the commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell while keeping it text-typed
# (matches the source file, where every data cell is an inline/shared string,
# even when the text happens to look numeric, e.g. "566.56" or "0.433").
function Set-TextCell($cell, $text) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextCell D2 '62.113.23'
Set-TextCell E2 '  +0.38%  '
Set-TextCell D3 '2.902.96'
Set-TextCell E3 '  -0.39%  '
Set-TextCell E4 '  +0.22%  '
Set-TextCell D5 '566.56'
Set-TextCell E5 '  -3.60%  '
Set-TextCell D6 '144.02'
Set-TextCell E6 '  -1.60%  '
Set-TextCell E7 '  +0.06%  '
Set-TextCell D8 '2.899.19'
Set-TextCell E8 '  -0.46%  '
Set-TextCell E9 '  -1.35%  '
Set-TextCell D10 '6.95'
Set-TextCell E10 '  -1.23%  '
Set-TextCell E11 '  -1.53%  '
Set-TextCell D12 '0.433'
Set-TextCell E12 '  -0.92%  '
Set-TextCell D13 '0.0000238'
Set-TextCell E13 '  -0.10%  '
Set-TextCell D14 '32.49'
Set-TextCell E14 '  -1.49%  '
Set-TextCell E15 '  +0.47%  '
Set-TextCell D16 '3.387.10'
Set-TextCell E16 '  -0.34%  '
Set-TextCell D17 '62.071.35'
Set-TextCell E17 '  +0.30%  '
Set-TextCell B18 'WrappedEther'
Set-TextCell C18 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell D18 '2.920.01'
Set-TextCell E18 '  +0.25%  '
Set-TextCell B19 'Polkadot'
Set-TextCell C19 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell D19 '6.58'
Set-TextCell E19 '  -0.54%  '
Set-TextCell D20 '429.27'
Set-TextCell E20 '  -1.52%  '
Set-TextCell D21 '13.05'
Set-TextCell E21 '  -3.05%  '
Set-TextCell D22 '0.658'
Set-TextCell E22 '  -0.27%  '
Set-TextCell D23 '6.84'
Set-TextCell E23 '  -1.55%  '
Set-TextCell E24 '  -3.13%  '
Set-TextCell E25 '  +0.66%  '
Set-TextCell B26 'Dai'
Set-TextCell C26 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell D26 '1.00'
Set-TextCell E26 '  -0.09%  '
Set-TextCell B27 'RenderToken'
Set-TextCell C27 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell D27 '10.00'
Set-TextCell E27 '  -2.33%  '
Set-TextCell E28 '  -3.91%  '
Set-TextCell D29 '0.0000111'
Set-TextCell E29 '  +3.68%  '
Set-TextCell E30 '  -4.75%  '
Set-TextCell E31 '  -2.64%  '
Set-TextCell E32 '  -5.02%  '
Set-TextCell E33 '  +0.05%  '
Set-TextCell B34 'Hedera'
Set-TextCell C34 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell D34 '0.106'
Set-TextCell E34 '  -3.45%  '
Set-TextCell B35 'EthereumClassic'
Set-TextCell C35 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell D35 '25.65'
Set-TextCell E35 '  -1.17%  '
Set-TextCell D36 '0.956'
Set-TextCell E36 '  -1.91%  '
Set-TextCell D37 '5.37'
Set-TextCell E37 '  -2.77%  '
Set-TextCell D38 '2.94'
Set-TextCell E38 '  -4.21%  '
Set-TextCell D39 '48.82'
Set-TextCell E39 '  -0.61%  '
Set-TextCell E40 '  -6.00%  '
Set-TextCell E41 '  -1.87%  '
Set-TextCell E42 '  -3.08%  '
Set-TextCell D43 '40.78'
Set-TextCell E43 '  +4.54%  '
Set-TextCell B44 'Maker'
Set-TextCell C44 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell D44 '2.712.74'
Set-TextCell E44 '  +0.47%  '
Set-TextCell B45 'TheGraph'
Set-TextCell C45 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell D45 '0.265'
Set-TextCell E45 '  -2.89%  '
Set-TextCell D46 '133.29'
Set-TextCell E46 '  -0.87%  '
Set-TextCell D47 '0.0335'
Set-TextCell E47 '  -0.88%  '
Set-TextCell D48 '346.85'
Set-TextCell E48 '  +0.05%  '
Set-TextCell E49 '  -0.04%  '
Set-TextCell D50 '0.000219'
Set-TextCell E50 '  +13.20%  '
Set-TextCell E51 '  -1.19%  '
